# Apply updated crypto price/volume data (GitHub Actions scrape refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''26.556.48'
$ws.Range('E2').Value = '  -0.12%  '
$ws.Range('D3').Value = '''1.852.71'
$ws.Range('E3').Value = '  +0.05%  '
$ws.Range('D4').Value = '''0.9994'
$ws.Range('E4').Value = '  -0.38%  '
$ws.Range('D5').Value = '''265.46'
$ws.Range('E5').Value = '  +2.09%  '
$ws.Range('D6').Value = '''0.9991'
$ws.Range('E6').Value = '  -0.43%  '
$ws.Range('D7').Value = '''0.5231'
$ws.Range('E7').Value = '  -0.39%  '
$ws.Range('E8').Value = '  -0.12%  '
$ws.Range('D9').Value = '''0.06816'
$ws.Range('E9').Value = '  +1.06%  '
$ws.Range('D10').Value = '''18.88'
$ws.Range('E10').Value = '  -2.95%  '
$ws.Range('D11').Value = '''0.7800'
$ws.Range('E11').Value = '  +0.44%  '
$ws.Range('D12').Value = '''0.07783'
$ws.Range('E12').Value = '  +1.00%  '
$ws.Range('D13').Value = '''1.861.85'
$ws.Range('E13').Value = '  -0.08%  '
$ws.Range('D14').Value = '''88.54'
$ws.Range('E14').Value = '  -0.57%  '
$ws.Range('E15').Value = '  -0.61%  '
$ws.Range('D16').Value = '''0.9994'
$ws.Range('E16').Value = '  -0.27%  '
$ws.Range('D17').Value = '''13.99'
$ws.Range('E17').Value = '  -1.66%  '
$ws.Range('D18').Value = '''0.000007976'
$ws.Range('E18').Value = '  +0.90%  '
$ws.Range('D19').Value = '''0.9985'
$ws.Range('E19').Value = '  -0.51%  '
$ws.Range('D20').Value = '''26.588.03'
$ws.Range('E20').Value = '  -0.12%  '
$ws.Range('D21').Value = '''2.090.69'
$ws.Range('E21').Value = '  +0.93%  '
$ws.Range('D22').Value = '''4.651'
$ws.Range('E22').Value = '  +0.61%  '
$ws.Range('D23').Value = '''9.567'
$ws.Range('E23').Value = '  -1.71%  '
$ws.Range('D24').Value = '''5.996'
$ws.Range('E24').Value = '  -0.08%  '
$ws.Range('D25').Value = '''144.79'
$ws.Range('E25').Value = '  +0.07%  '
$ws.Range('D26').Value = '''2.231'
$ws.Range('E26').Value = '  -6.02%  '
$ws.Range('D27').Value = '''1.660'
$ws.Range('E27').Value = '  -0.19%  '
$ws.Range('D28').Value = '''17.05'
$ws.Range('E28').Value = '  -0.18%  '
$ws.Range('D29').Value = '''112.39'
$ws.Range('E29').Value = '  +1.02%  '
$ws.Range('D30').Value = '''4.204'
$ws.Range('E30').Value = '  -0.65%  '
$ws.Range('D31').Value = '''4.147'
$ws.Range('E31').Value = '  -1.24%  '
$ws.Range('D32').Value = '''0.08763'
$ws.Range('E32').Value = '  -0.37%  '
$ws.Range('D33').Value = '''0.04854'
$ws.Range('E33').Value = '  -0.69%  '
$ws.Range('D34').Value = '''1.141'
$ws.Range('E34').Value = '  -0.39%  '
$ws.Range('D35').Value = '''0.7210'
$ws.Range('E35').Value = '  +1.31%  '
$ws.Range('D36').Value = '''2.853'
$ws.Range('E36').Value = '  -0.97%  '
$ws.Range('D37').Value = '''3.109'
$ws.Range('E37').Value = '  -1.09%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = '''0.01788'
$ws.Range('E38').Value = '  -1.69%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').Value = '''2.227'
$ws.Range('E39').Value = '  +1.04%  '
$ws.Range('D40').Value = '''0.4900'
$ws.Range('E40').Value = '  -1.70%  '
$ws.Range('D41').Value = '''0.9189'
$ws.Range('E41').Value = '  +1.91%  '
$ws.Range('D42').Value = '''111.82'
$ws.Range('E42').Value = '  -2.72%  '
$ws.Range('D43').Value = '''6.081'
$ws.Range('E43').Value = '  -0.06%  '
$ws.Range('D44').Value = '''0.9983'
$ws.Range('E44').Value = '  -0.59%  '
$ws.Range('D45').Value = '''7.759'
$ws.Range('E45').Value = '  -0.59%  '
$ws.Range('D46').Value = '''0.4201'
$ws.Range('E46').Value = '  -2.61%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = '''9.177'
$ws.Range('E47').Value = '  +0.47%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').Value = '''0.05952'
$ws.Range('E48').Value = '  +0.48%  '
$ws.Range('D49').Value = '''0.1245'
$ws.Range('E49').Value = '  -3.91%  '
$ws.Range('D50').Value = '''35.13'
$ws.Range('E50').Value = '  -1.20%  '
$ws.Range('D51').Value = '''0.8912'
$ws.Range('E51').Value = '  +2.71%  '
